$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows of compound data (rows 2-9, columns A-X) are being reordered.
# New order (top to bottom) of compounds for rows 2..9:
#   row2 -> naphthalene   (was row5)
#   row3 -> dodecane      (was row2)
#   row4 -> phenol        (unchanged, stays row4)
#   row5 -> oleic acid    (was row8)
#   row6 -> palmitic acid (was row3)
#   row7 -> p-dichlorobenzene / 1,4-dichlorobenzene (was row6, name in col A changes)
#   row8 -> notvalidcomp / unidentified (was row7)
#   row9 -> capric acid   (unchanged, stays row9)

$data = @{
    2 = @("naphthalene", "naphthalene", "C10H8", "C1=CC=C2C=CC=CC2=C1", 128.17, 3.3, 10, 0, 8, 0, 0.9371147694468284, 0, 0.06291643910431459, 0, 0, 10, 0, 0, 0, 0, 1.000031208551143, 0, 0, 0)
    3 = @("dodecane", "dodecane", "C12H26", "CCCCCCCCCCCC", 170.33, 6.1, 12, 0, 26, 0, 0.846192684788352, 0, 0.1538660247754359, 0, 12, 0, 0, 0, 0, 1.000058709563788, 0, 0, 0, 0)
    5 = @("oleic acid", "(z)-octadec-9-enoic acid", "C18H34O2", "CCCCCCCCC=CCCCCCCCC(=O)O", 282.5, 6.5, 18, 0, 34, 2, 0.7653026548672566, 0, 0.121316814159292, 0.1132672566371681, 17, 0, 0, 0, 1, 0.8405345132743363, 0, 0, 0, 0.1593522123893805)
    6 = @("palmitic acid", "hexadecanoic acid", "C16H32O2", "CCCCCCCCCCCCCCCC(=O)O", 256.42, 6.4, 16, 0, 32, 2, 0.7494579205990172, 0, 0.125793619842446, 0.1247874580765931, 15, 0, 0, 0, 1, 0.8244793697839481, 0, 0, 0, 0.1755596287341081)
    7 = @("p-dichlorobenzene", "1,4-dichlorobenzene", "C6H4Cl2", "C1=CC(=CC=C1Cl)Cl", 147, 3.4, 6, 2, 4, 0, 0.4902448979591837, 0.4823129251700681, 0.02742857142857143, 0, 0, 6, 2, 0, 0, 0, 0.517673469387755, 0.4823129251700681, 0, 0)
    8 = @("notvalidcomp", "unidentified")
}

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X")

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $col = $cols[$i]
        $addr = "$col$r"
        if ($i -lt $rowVals.Length) {
            $ws.Range($addr).Value = $rowVals[$i]
        } else {
            $ws.Range($addr).Value = ""
        }
    }
}
